$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header text cells (Volume/Number + report week dates)
# ---------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/13/2025  Through  10/19/2025"


# --- Simple numeric value changes ---
$ws.Range("M14").Value = -25
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("M15").Value = 27.272727272727
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -18.181818181818
$ws.Range("I16").Value = 86
$ws.Range("J16").Value = 95
$ws.Range("K16").Value = -9.473684210526
$ws.Range("L16").Value = -25.217391304347
$ws.Range("M16").Value = -58.454106280193
$ws.Range("N16").Value = -88.684210526315
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -39.130434782608
$ws.Range("I17").Value = 182
$ws.Range("J17").Value = 235
$ws.Range("K17").Value = -22.553191489361
$ws.Range("L17").Value = -3.703703703703
$ws.Range("M17").Value = 124.691358024691
$ws.Range("N17").Value = -8.080808080808
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 14
$ws.Range("H18").Value = -17.647058823529
$ws.Range("I18").Value = 164
$ws.Range("J18").Value = 178
$ws.Range("K18").Value = -7.865168539325
$ws.Range("L18").Value = -23.720930232558
$ws.Range("M18").Value = -27.433628318584
$ws.Range("N18").Value = -85.922746781115
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -46.153846153846
$ws.Range("F19").Value = 51
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = 15.90909090909
$ws.Range("I19").Value = 462
$ws.Range("J19").Value = 457
$ws.Range("K19").Value = 1.094091903719
$ws.Range("L19").Value = -5.714285714285
$ws.Range("M19").Value = 33.91304347826
$ws.Range("N19").Value = -6.097560975609
$ws.Range("C20").Value = 16
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 77.777777777777
$ws.Range("F20").Value = 34
$ws.Range("G20").Value = 42
$ws.Range("H20").Value = -19.047619047619
$ws.Range("I20").Value = 282
$ws.Range("J20").Value = 313
$ws.Range("K20").Value = -9.904153354632
$ws.Range("L20").Value = -1.742160278745
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -92.885973763874
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = 2.857142857142
$ws.Range("F21").Value = 124
$ws.Range("G21").Value = 141
$ws.Range("H21").Value = -12.056737588652
$ws.Range("I21").Value = 1193
$ws.Range("J21").Value = 1300
$ws.Range("K21").Value = -8.230769230769
$ws.Range("L21").Value = -8.931297709923
$ws.Range("M21").Value = 12.335216572504
$ws.Range("N21").Value = -81.951588502269
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("J22").Value = 29
$ws.Range("K22").Value = -27.586206896551
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("M23").Value = 62.962962962963
$ws.Range("C24").Value = 20
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -23.076923076923
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = -5.050505050505
$ws.Range("I24").Value = 878
$ws.Range("J24").Value = 967
$ws.Range("K24").Value = -9.203722854188
$ws.Range("L24").Value = -23.784722222222
$ws.Range("M24").Value = 7.995079950799
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -86.666666666666
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -30.952380952381
$ws.Range("I25").Value = 291
$ws.Range("J25").Value = 415
$ws.Range("K25").Value = -29.879518072289
$ws.Range("L25").Value = -30.714285714285
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -37.5
$ws.Range("F26").Value = 37
$ws.Range("G26").Value = 50
$ws.Range("H26").Value = -26
$ws.Range("I26").Value = 424
$ws.Range("J26").Value = 424
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 9.560723514211
$ws.Range("M26").Value = 31.269349845201
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("E28").Value = -100
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("J28").Value = 33
$ws.Range("K28").Value = 42.424242424242
$ws.Range("L28").Value = 27.027027027027
$ws.Range("M29").Value = 0
$ws.Range("M30").Value = -20

# --- Style-changing cells (value + type + number format) ---
# Text/N/A template (style 13): C22  |  Count template (style 14): F15  |  Pct template (style 15): K15
$ws.Range("C15").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("K15").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("C28").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D31").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E31").PasteSpecial(-4122)
